$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header labels:
#   columns 1-10  (A-J)  "<name>_old" -> "<name>_FV2210"
#   columns 12-21 (L-U)  "<name>_new" -> "<name>_FV2304"
for ($i = 1; $i -le 10; $i++) {
    $cell = $ws.Cells.Item(1, $i)
    $label = $cell.Value()
    $cell.Value = $label.Replace("_old", "_FV2210")
}
for ($i = 12; $i -le 21; $i++) {
    $cell = $ws.Cells.Item(1, $i)
    $label = $cell.Value()
    $cell.Value = $label.Replace("_new", "_FV2304")
}

# Turn the used range into an Excel table ("Table1") with an autofilter.
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U67"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
